$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "Бронь закрыта" (booking closed) column is being inserted as column B.
# The former column B ("id спектакль") shifts to C, and the former column C
# ("id ценовая политика") shifts to D. Read the old values first, then write
# them into their new homes (reading a COM property requires calling it,
# i.e. `.Value()`, in this runtime).
$oldC1 = $ws.Range("C1").Value()
$oldC2 = $ws.Range("C2").Value()
$oldC3 = $ws.Range("C3").Value()

$oldB1 = $ws.Range("B1").Value()
$oldB2 = $ws.Range("B2").Value()
$oldB3 = $ws.Range("B3").Value()

# Move old column C ("id ценовая политика") into new column D.
$ws.Range("D1").Value = $oldC1
$ws.Range("D2").Value = $oldC2
$ws.Range("D3").Value = $oldC3

# Move old column B ("id спектакль") into new column C.
$ws.Range("C1").Value = $oldB1
$ws.Range("C2").Value = $oldB2
$ws.Range("C3").Value = $oldB3

# Populate the new "Бронь закрыта" column B with Д (yes) / Н (no) markers.
$ws.Range("B1").Value = "Бронь закрыта"
$ws.Range("B2").Value = "Д"
$ws.Range("B3").Value = "Н"

# Match the author's final selection (whole used range, anchored at A1).
$ws.Range("A1:D3").Select()
